# zemoso interview missed questions
$wb = $excel.ActiveWorkbook

# Add the two missed interview questions to their respective sheets.

# "nodejs" sheet: append a new question about fork vs spawn right after
# the existing last row (row 26 -> new row 27).
$wsNode = $wb.Worksheets.Item("nodejs")
$wsNode.Range("A27").Value = "Difference between fork and spwan"

# "Design Patterns" sheet: append a new question about the CAP theorem
# right after the existing last row (row 4 -> new row 5).
$wsDesign = $wb.Worksheets.Item("Design Patterns")
$wsDesign.Range("A5").Value = "CAP Therome"

# Update selections/active cells to reflect where the user left off editing.
$wsNode.Range("A28").Select()

$wsDesign.Activate()
$wsDesign.Range("A5").Select()
